# Update "countries & provincias Spain" data dump with newer figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 01:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 465689
$ws.Range("C4").Value = 30659
$ws.Range("D4").Value = 25156
$ws.Range("E4").Value = 423985
$ws.Range("G4").Value = 1757
$ws.Range("H4").Value = 16548

# Row 7 - Alemania
$ws.Range("B7").Value = 118235
$ws.Range("C7").Value = 4939
$ws.Range("E7").Value = 63221

# Row 16 - Canada
$ws.Range("B16").Value = 20765
$ws.Range("C16").Value = 1327
$ws.Range("D16").Value = 5311
$ws.Range("E16").Value = 14945
$ws.Range("G16").Value = 82
$ws.Range("H16").Value = 509

# Row 31 - Chequia
$ws.Range("B31").Value = 5569
$ws.Range("C31").Value = 257
$ws.Range("E31").Value = 5156

# Row 53 - Singapur
$ws.Range("E53").Value = 1443
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 7

# Row 55 - Argentina
$ws.Range("E55").Value = 1358
$ws.Range("G55").Value = 7
$ws.Range("H55").Value = 72

# Row 97 - Ghana
$ws.Range("D97").Value = 3
$ws.Range("E97").Value = 369

# Row 117 - Mayotte
$ws.Range("D117").Value = 26
$ws.Range("E117").Value = 156
$ws.Range("F117").Value = 4

# Row 122 - Guadalupe
$ws.Range("B122").Value = 143
$ws.Range("C122").Value = 2
$ws.Range("D122").Value = 67
$ws.Range("E122").Value = 68

# Row 123 - Republica de Yibuti
$ws.Range("E123").Value = 109
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 1

# Row 150 - Bahamas
$ws.Range("B150").Value = 41
$ws.Range("C150").Value = 1
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 8
